# Apply the data-refresh edit to the bread_coop crawl output workbook.
#
# The site was re-crawled later the same day (crawl timestamp moved from
# 2022-12-27 06:49:19 to 2022-12-27 12:55:54 for every data row), and a
# handful of products show updated rating counts / values, while a few
# "Online kein Bestand" (out of stock) labels toggled on/off in the
# productAriaLabel (column M) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-12-27 06:49:19"
$newTimestamp = "2022-12-27 12:55:54"

# Column indexes (1-based): A=1 id, B=2 title, C=3 href, D=4 ratingAmount,
# E=5 ratingValue, ... M=13 productAriaLabel, N=14 declarationIcons,
# O=15 timestamp.
$colD = 4
$colE = 5
$colM = 13
$colO = 15

$lastRow = 412

# ratingAmount (and occasionally ratingValue) changes, keyed by row number.
$ratingAmountChanges = @{
    11  = 24
    18  = 21
    19  = 46
    23  = 22
    29  = 5
    31  = 12
    36  = 5
    37  = 14
    43  = 6
    48  = 34
    49  = 20
    53  = 10
    54  = 14
    55  = 6
    75  = 22
    114 = 6
    120 = 42
    127 = 13
    138 = 5
    147 = 6
    150 = 12
    203 = 5
    217 = 3
    254 = 33
}

$ratingValueChanges = @{
    147 = 4.5
    254 = 3.5
}

foreach ($row in $ratingAmountChanges.Keys) {
    $ws.Cells.Item($row, $colD).Value = $ratingAmountChanges[$row]
}

foreach ($row in $ratingValueChanges.Keys) {
    $ws.Cells.Item($row, $colE).Value = $ratingValueChanges[$row]
}

# productAriaLabel (column M) text updates - stock-status wording changed.
$ws.Cells.Item(266, $colM).Value = "St Michel Madeleines 3.30 Schweizer Franken"
$ws.Cells.Item(273, $colM).Value = "St Michel Madeleines 10 Stück 2.80 Schweizer Franken"
$ws.Cells.Item(324, $colM).Value = "Prix Garantie Madeleines choco 2.80 Schweizer Franken"
$ws.Cells.Item(404, $colM).Value = "Buitoni Pizzateig Glutenfrei &amp; Ohne Lactose Rund Ausgewallt Ø25cm - Online kein Bestand 4.95 Schweizer Franken"

# Every data row (2..412) got a refreshed crawl timestamp in column O.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colO).Value = $newTimestamp
}
